# Added Configurable zero_before_threshold parameter to enable setting dims
# before noise_threshold or First Rise Point to 0.
#
# This updates the First_Noticeable_Increase_Index (C), the
# First_Noticeable_Increase_Cumulative_Value (E), and the Pulse_Width (G)
# columns on each of the four "Step3_DataPts_*" sheets to reflect the new
# zero_before_threshold behavior.

$wb = $excel.ActiveWorkbook

# New values for First_Noticeable_Increase_Index (column C) and
# First_Noticeable_Increase_Cumulative_Value (column E) - these are the
# same across all four Step3_DataPts_* sheets for a given row.
$cValues = @{ 2 = 87; 3 = 87; 4 = 87; 5 = 88; 6 = 88 }
$eValues = @{
    2 = 0.002765405891476024
    3 = 0.01635860265877747
    4 = 0.0214549997476814
    5 = 0.0101733135263407
    6 = 0.004015533523784611
}

# New Pulse_Width (column G) values, specific to each threshold sheet.
$gValuesBySheet = @{
    "Step3_DataPts_0.5" = @{ 2 = 21; 3 = 22; 4 = 21; 5 = 22; 6 = 20 }
    "Step3_DataPts_0.7" = @{ 2 = 45; 3 = 53; 4 = 47; 5 = 54; 6 = 45 }
    "Step3_DataPts_0.8" = @{ 2 = 70; 3 = 69; 4 = 68; 5 = 67; 6 = 69 }
    "Step3_DataPts_0.9" = @{ 2 = 79; 3 = 80; 4 = 79; 5 = 78; 6 = 80 }
}

foreach ($sheetName in $gValuesBySheet.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $gValues = $gValuesBySheet[$sheetName]

    foreach ($row in 2..6) {
        $ws.Cells.Item($row, 3).Value = $cValues[$row]
        $ws.Cells.Item($row, 5).Value = $eValues[$row]
        $ws.Cells.Item($row, 7).Value = $gValues[$row]
    }
}
